# chore: update Sheets via scheduled runner
#
# Refreshes the cached market-board figures (currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ, columns H:N) for a batch of leve rows across
# the profession sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 122.8
$ws.Range("I4").Value = 122.8
$ws.Range("K4").Value = 122.8
$ws.Range("M4").Value = -8.799999999999997

$ws.Range("H6").Value = 217.375
$ws.Range("I6").Value = 177.8
$ws.Range("K6").Value = 533.4000000000001
$ws.Range("M6").Value = -421.4000000000001

$ws.Range("H9").Value = 200
$ws.Range("J9").Value = 200
$ws.Range("L9").Value = 200
$ws.Range("N9").Value = -538

$ws.Range("H31").Value = 552
$ws.Range("I31").Value = 100
$ws.Range("J31").Value = 1004
$ws.Range("K31").Value = 300
$ws.Range("L31").Value = 3012
$ws.Range("M31").Value = -70
$ws.Range("N31").Value = -3472

$ws.Range("H32").Value = 875.46155
$ws.Range("I32").Value = 823
$ws.Range("K32").Value = 823
$ws.Range("M32").Value = -497

$ws.Range("H38").Value = 695.8
$ws.Range("I38").Value = 57.25
$ws.Range("K38").Value = 171.75
$ws.Range("M38").Value = 200.25

$ws.Range("H39").Value = 352.625
$ws.Range("I39").Value = 144.17647
$ws.Range("J39").Value = 858.8570999999999
$ws.Range("K39").Value = 432.52941
$ws.Range("L39").Value = 2576.5713
$ws.Range("M39").Value = -136.52941
$ws.Range("N39").Value = -3168.5713

$ws.Range("H132").Value = 2932.5881
$ws.Range("I132").Value = 2932.5881
$ws.Range("K132").Value = 8797.764299999999
$ws.Range("M132").Value = -6267.764299999999

$ws.Range("H137").Value = 2670.1428
$ws.Range("I137").Value = 2372.75
$ws.Range("K137").Value = 7118.25
$ws.Range("M137").Value = -4568.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 999.5
$ws.Range("I102").Value = 999.5
$ws.Range("K102").Value = 999.5
$ws.Range("M102").Value = 622.5

$ws.Range("H122").Value = 3500
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -19900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1500
$ws.Range("I20").Value = 1500
$ws.Range("K20").Value = 1500
$ws.Range("M20").Value = -1253

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5483.0835
$ws.Range("I31").Value = 4413.857
$ws.Range("J31").Value = 6980
$ws.Range("K31").Value = 4413.857
$ws.Range("L31").Value = 6980
$ws.Range("M31").Value = -4118.857
$ws.Range("N31").Value = -7570

$ws.Range("H34").Value = 5483.0835
$ws.Range("I34").Value = 4413.857
$ws.Range("J34").Value = 6980
$ws.Range("K34").Value = 4413.857
$ws.Range("L34").Value = 6980
$ws.Range("M34").Value = -4211.857
$ws.Range("N34").Value = -7384

$ws.Range("H36").Value = 13000
$ws.Range("I36").Value = 14000
$ws.Range("J36").Value = 12000
$ws.Range("K36").Value = 14000
$ws.Range("L36").Value = 12000
$ws.Range("M36").Value = -13612
$ws.Range("N36").Value = -12776

$ws.Range("H40").Value = 13000
$ws.Range("I40").Value = 14000
$ws.Range("J40").Value = 12000
$ws.Range("K40").Value = 14000
$ws.Range("L40").Value = 12000
$ws.Range("M40").Value = -13840
$ws.Range("N40").Value = -12320

$ws.Range("H57").Value = 20000
$ws.Range("J57").Value = 20000
$ws.Range("L57").Value = 20000
$ws.Range("N57").Value = -21120

$ws.Range("H62").Value = 5899.75
$ws.Range("I62").Value = 5866.3335
$ws.Range("K62").Value = 5866.3335
$ws.Range("M62").Value = -5242.3335

$ws.Range("H65").Value = 5899.75
$ws.Range("I65").Value = 5866.3335
$ws.Range("K65").Value = 29331.6675
$ws.Range("M65").Value = -26211.6675

$ws.Range("H122").Value = 1486.5555
$ws.Range("I122").Value = 1141.6666
$ws.Range("J122").Value = 1659
$ws.Range("K122").Value = 3424.9998
$ws.Range("L122").Value = 4977
$ws.Range("M122").Value = -974.9998000000001
$ws.Range("N122").Value = -9877

$ws.Range("H132").Value = 7331
$ws.Range("J132").Value = 6994
$ws.Range("L132").Value = 20982
$ws.Range("N132").Value = -26042

$ws.Range("H134").Value = 3431.8635
$ws.Range("I134").Value = 3506.85
$ws.Range("J134").Value = 2682
$ws.Range("K134").Value = 10520.55
$ws.Range("L134").Value = 8046
$ws.Range("M134").Value = -7985.549999999999
$ws.Range("N134").Value = -13116

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 403.0909
$ws.Range("I7").Value = 195.25
$ws.Range("K7").Value = 585.75
$ws.Range("M7").Value = -473.75

$ws.Range("H32").Value = 875
$ws.Range("J32").Value = 1000
$ws.Range("L32").Value = 3000
$ws.Range("N32").Value = -3566

$ws.Range("H51").Value = 999.5
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H92").Value = 1460
$ws.Range("I92").Value = 1625
$ws.Range("K92").Value = 4875
$ws.Range("M92").Value = -3627

$ws.Range("H99").Value = 2431.25
$ws.Range("I99").Value = 2431.25
$ws.Range("K99").Value = 7293.75
$ws.Range("M99").Value = -5047.75

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 30001
$ws.Range("J47").Value = 30001
$ws.Range("L47").Value = 30001
$ws.Range("N47").Value = -31137

$ws.Range("H55").Value = 15999.889
$ws.Range("I55").Value = 5200
$ws.Range("J55").Value = 24639.8
$ws.Range("K55").Value = 5200
$ws.Range("L55").Value = 24639.8
$ws.Range("M55").Value = -4873
$ws.Range("N55").Value = -25293.8

$ws.Range("H122").Value = 5001893.5
$ws.Range("I122").Value = 5001893.5
$ws.Range("K122").Value = 15005680.5
$ws.Range("M122").Value = -15003230.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3336
$ws.Range("I46").Value = 3557.1428
$ws.Range("K46").Value = 3557.1428
$ws.Range("M46").Value = -3369.1428

$ws.Range("H61").Value = 6166.6665
$ws.Range("I61").Value = 5000
$ws.Range("K61").Value = 5000
$ws.Range("M61").Value = -4798

$ws.Range("H113").Value = 6166.6665
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -2830

$ws.Range("H132").Value = 12728.429
$ws.Range("I132").Value = 4819.8
$ws.Range("J132").Value = 32500
$ws.Range("K132").Value = 14459.4
$ws.Range("L132").Value = 97500
$ws.Range("M132").Value = -11929.4
$ws.Range("N132").Value = -102560

$ws.Range("H136").Value = 3301.875
$ws.Range("I136").Value = 3301.875
$ws.Range("K136").Value = 9905.625
$ws.Range("M136").Value = -7355.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 37032.5
$ws.Range("I48").Value = 25000
$ws.Range("J48").Value = 49065
$ws.Range("K48").Value = 25000
$ws.Range("L48").Value = 49065
$ws.Range("M48").Value = -24431
$ws.Range("N48").Value = -50203

$ws.Range("H122").Value = 1099.5
$ws.Range("I122").Value = 1099.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3298.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -848.5
$ws.Range("N122").ClearContents()

$ws.Range("H135").Value = 40000
$ws.Range("J135").Value = 40000
$ws.Range("L135").Value = 40000
$ws.Range("N135").Value = -50140

$ws.Range("H136").Value = 1551
$ws.Range("I136").Value = 1006.7778
$ws.Range("K136").Value = 3020.3334
$ws.Range("M136").Value = -470.3334
